$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo in the "Years of Experience" header (column B)
$ws.Range("B1").Value = "Years of Experience"

# Apply header styling: bold font, orange fill, thin border around each cell
$headerRange = $ws.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.Interior.Color = 42495
$headerRange.Borders.LineStyle = 1
